{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" answer strings found in the\n// table cells of the document with their new values, one-for-one.\n// Each old value is unique in the document, so a simple exact search +\n// full replace of the matched range is safe and keeps the surrounding\n// run formatting (font, size, etc.) intact.\nconst replacements = [\n  [\"287\u00f78=35, 7\", \"818\u00f77=116, 6\"],\n  [\"713\u00f79=79, 2\", \"440\u00f75=88, 0\"],\n  [\"644\u00f72=322, 0\", \"737\u00f73=245, 2\"],\n  [\"122\u00f73=40, 2\", \"545\u00f75=109, 0\"],\n  [\"387\u00f79=43, 0\", \"358\u00f76=59, 4\"],\n  [\"432\u00f77=61, 5\", \"427\u00f73=142, 1\"],\n  [\"464\u00f76=77, 2\", \"104\u00f77=14, 6\"],\n  [\"513\u00f73=171, 0\", \"402\u00f75=80, 2\"],\n  [\"156\u00f79=17, 3\", \"557\u00f78=69, 5\"],\n  [\"425\u00f74=106, 1\", \"604\u00f78=75, 4\"],\n  [\"684\u00f79=76, 0\", \"107\u00f72=53, 1\"],\n  [\"749\u00f72=374, 1\", \"482\u00f75=96, 2\"],\n  [\"228\u00f72=114, 0\", \"540\u00f73=180, 0\"],\n  [\"996\u00f79=110, 6\", \"382\u00f76=63, 4\"],\n  [\"371\u00f79=41, 2\", \"628\u00f75=125, 3\"],\n  [\"602\u00f73=200, 2\", \"846\u00f73=282, 0\"],\n  [\"522\u00f75=104, 2\", \"978\u00f72=489, 0\"],\n  [\"127\u00f72=63, 1\", \"183\u00f78=22, 7\"],\n  [\"744\u00f75=148, 4\", \"198\u00f72=99, 0\"],\n  [\"769\u00f75=153, 4\", \"364\u00f77=52, 0\"],\n  [\"120\u00f73=40, 0\", \"824\u00f77=117, 5\"],\n  [\"984\u00f72=492, 0\", \"938\u00f75=187, 3\"],\n  [\"560\u00f75=112, 0\", \"502\u00f74=125, 2\"],\n  [\"336\u00f72=168, 0\", \"609\u00f79=67, 6\"],\n  [\"582\u00f76=97, 0\", \"917\u00f77=131, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit / one-digit\" division-answer strings in the\n# table cells with their new values, one-for-one. Each old value is unique\n# in the document, so Find/Replace (exact match, whole document) for each\n# pair is safe and preserves the surrounding run formatting.\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"287\u00f78=35, 7\", \"818\u00f77=116, 6\"),\n    @(\"713\u00f79=79, 2\", \"440\u00f75=88, 0\"),\n    @(\"644\u00f72=322, 0\", \"737\u00f73=245, 2\"),\n    @(\"122\u00f73=40, 2\", \"545\u00f75=109, 0\"),\n    @(\"387\u00f79=43, 0\", \"358\u00f76=59, 4\"),\n    @(\"432\u00f77=61, 5\", \"427\u00f73=142, 1\"),\n    @(\"464\u00f76=77, 2\", \"104\u00f77=14, 6\"),\n    @(\"513\u00f73=171, 0\", \"402\u00f75=80, 2\"),\n    @(\"156\u00f79=17, 3\", \"557\u00f78=69, 5\"),\n    @(\"425\u00f74=106, 1\", \"604\u00f78=75, 4\"),\n    @(\"684\u00f79=76, 0\", \"107\u00f72=53, 1\"),\n    @(\"749\u00f72=374, 1\", \"482\u00f75=96, 2\"),\n    @(\"228\u00f72=114, 0\", \"540\u00f73=180, 0\"),\n    @(\"996\u00f79=110, 6\", \"382\u00f76=63, 4\"),\n    @(\"371\u00f79=41, 2\", \"628\u00f75=125, 3\"),\n    @(\"602\u00f73=200, 2\", \"846\u00f73=282, 0\"),\n    @(\"522\u00f75=104, 2\", \"978\u00f72=489, 0\"),\n    @(\"127\u00f72=63, 1\", \"183\u00f78=22, 7\"),\n    @(\"744\u00f75=148, 4\", \"198\u00f72=99, 0\"),\n    @(\"769\u00f75=153, 4\", \"364\u00f77=52, 0\"),\n    @(\"120\u00f73=40, 0\", \"824\u00f77=117, 5\"),\n    @(\"984\u00f72=492, 0\", \"938\u00f75=187, 3\"),\n    @(\"560\u00f75=112, 0\", \"502\u00f74=125, 2\"),\n    @(\"336\u00f72=168, 0\", \"609\u00f79=67, 6\"),\n    @(\"582\u00f76=97, 0\", \"917\u00f77=131, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n\nWrite-Output \"done\"\n"}
